$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 431.13043
$ws.Range("I33").Value = 435.4091
$ws.Range("K33").Value = 435.4091
$ws.Range("M33").Value = -206.4091
$ws.Range("H40").Value = 3698.9333
$ws.Range("J40").Value = 5360
$ws.Range("L40").Value = 5360
$ws.Range("N40").Value = -5710
$ws.Range("H80").Value = 594.7895
$ws.Range("I80").Value = 365.6
$ws.Range("K80").Value = 1096.8
$ws.Range("M80").Value = -98.80000000000018
$ws.Range("H83").Value = 594.7895
$ws.Range("I83").Value = 365.6
$ws.Range("K83").Value = 3290.4
$ws.Range("M83").Value = 1701.6
$ws.Range("H86").Value = 1799.7273
$ws.Range("I86").Value = 2160.6
$ws.Range("J86").Value = 1499
$ws.Range("K86").Value = 2160.6
$ws.Range("L86").Value = 1499
$ws.Range("M86").Value = -1037.6
$ws.Range("N86").Value = -3745
$ws.Range("H89").Value = 1799.7273
$ws.Range("I89").Value = 2160.6
$ws.Range("J89").Value = 1499
$ws.Range("K89").Value = 10803
$ws.Range("L89").Value = 7495
$ws.Range("M89").Value = -5187
$ws.Range("N89").Value = -18727
$ws.Range("H92").Value = 5001.3335
$ws.Range("I92").Value = 7002
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 7002
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = -5754
$ws.Range("N92").Value = -3496
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459
$ws.Range("H115").Value = 401.25
$ws.Range("I115").Value = 401.25
$ws.Range("K115").Value = 1203.75
$ws.Range("M115").Value = 363.25
$ws.Range("H116").Value = 4696
$ws.Range("I116").Value = 4169
$ws.Range("K116").Value = 4169
$ws.Range("M116").Value = -727
$ws.Range("H127").Value = 1654
$ws.Range("I127").Value = 1055
$ws.Range("K127").Value = 3165
$ws.Range("M127").Value = 1795
$ws.Range("H132").Value = 32261530
$ws.Range("I132").Value = 50004236
$ws.Range("J132").Value = 2064.7273
$ws.Range("K132").Value = 150012708
$ws.Range("L132").Value = 6194.1819
$ws.Range("M132").Value = -150010178
$ws.Range("N132").Value = -11254.1819
$ws.Range("H135").Value = 786.76746
$ws.Range("I135").Value = 766.91895
$ws.Range("K135").Value = 6902.27055
$ws.Range("M135").Value = -4367.27055
$ws.Range("H141").Value = 2016.875
$ws.Range("I141").Value = 2016.875
$ws.Range("K141").Value = 6050.625
$ws.Range("M141").Value = -870.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1024.6666
$ws.Range("J4").Value = 75
$ws.Range("L4").Value = 75
$ws.Range("N4").Value = -307
$ws.Range("H32").Value = 2693.8513
$ws.Range("I32").Value = 1406.4531
$ws.Range("J32").Value = 10933.2
$ws.Range("K32").Value = 1406.4531
$ws.Range("L32").Value = 10933.2
$ws.Range("M32").Value = -1119.4531
$ws.Range("N32").Value = -11507.2
$ws.Range("H97").Value = 1618738
$ws.Range("I97").Value = 1618738
$ws.Range("K97").Value = 1618738
$ws.Range("M97").Value = -1618242
$ws.Range("H102").Value = 3336462.2
$ws.Range("J102").Value = 5553.125
$ws.Range("L102").Value = 5553.125
$ws.Range("N102").Value = -8797.125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3666311.8
$ws.Range("I99").Value = 5717098
$ws.Range("J99").Value = 4193.0713
$ws.Range("K99").Value = 5717098
$ws.Range("L99").Value = 4193.0713
$ws.Range("M99").Value = -5715600
$ws.Range("N99").Value = -7189.0713
$ws.Range("H105").Value = 5683428
$ws.Range("I105").Value = 6945990
$ws.Range("K105").Value = 6945990
$ws.Range("M105").Value = -6944243
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2899.5
$ws.Range("I62").Value = 2849.5
$ws.Range("J62").Value = 2999.5
$ws.Range("K62").Value = 2849.5
$ws.Range("L62").Value = 2999.5
$ws.Range("M62").Value = -2225.5
$ws.Range("N62").Value = -4247.5
$ws.Range("H65").Value = 2899.5
$ws.Range("I65").Value = 2849.5
$ws.Range("J65").Value = 2999.5
$ws.Range("K65").Value = 14247.5
$ws.Range("L65").Value = 14997.5
$ws.Range("M65").Value = -11127.5
$ws.Range("N65").Value = -21237.5
$ws.Range("H105").Value = 3516.5
$ws.Range("I105").Value = 1805
$ws.Range("J105").Value = 4372.25
$ws.Range("K105").Value = 1805
$ws.Range("L105").Value = 4372.25
$ws.Range("M105").Value = -58
$ws.Range("N105").Value = -7866.25
$ws.Range("H107").Value = 936.8889
$ws.Range("I107").Value = 900.0741
$ws.Range("K107").Value = 900.0741
$ws.Range("M107").Value = 1019.9259
$ws.Range("H132").Value = 67209.55499999999
$ws.Range("I132").Value = 43705.918
$ws.Range("J132").Value = 180027
$ws.Range("K132").Value = 131117.754
$ws.Range("L132").Value = 540081
$ws.Range("M132").Value = -128587.754
$ws.Range("N132").Value = -545141
$ws.Range("H134").Value = 3135.318
$ws.Range("I134").Value = 1398.5834
$ws.Range("K134").Value = 4195.7502
$ws.Range("M134").Value = -1660.7502
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 55774.125
$ws.Range("J12").Value = 232.75
$ws.Range("L12").Value = 698.25
$ws.Range("N12").Value = -1044.25
$ws.Range("H86").Value = 394.8
$ws.Range("I86").Value = 316.33334
$ws.Range("K86").Value = 949.0000200000001
$ws.Range("M86").Value = 236.9999799999999
$ws.Range("H89").Value = 394.8
$ws.Range("I89").Value = 316.33334
$ws.Range("K89").Value = 2847.00006
$ws.Range("M89").Value = 3080.99994
$ws.Range("H132").Value = 1511.9375
$ws.Range("I132").Value = 1043.3334
$ws.Range("K132").Value = 9390.000599999999
$ws.Range("M132").Value = -6860.000599999999
$ws.Range("H136").Value = 1379.5555
$ws.Range("J136").Value = 3500
$ws.Range("L136").Value = 10500
$ws.Range("N136").Value = -20700
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6507.647
$ws.Range("I2").Value = 1048.1
$ws.Range("K2").Value = 1048.1
$ws.Range("M2").Value = -935.0999999999999
$ws.Range("H23").Value = 1376.75
$ws.Range("I23").Value = 992
$ws.Range("J23").Value = 1505
$ws.Range("K23").Value = 992
$ws.Range("L23").Value = 1505
$ws.Range("M23").Value = -769
$ws.Range("N23").Value = -1951
$ws.Range("H80").Value = 1745226.6
$ws.Range("J80").Value = 5781.2
$ws.Range("L80").Value = 5781.2
$ws.Range("N80").Value = -7777.2
$ws.Range("H83").Value = 1745226.6
$ws.Range("J83").Value = 5781.2
$ws.Range("L83").Value = 28906
$ws.Range("N83").Value = -38890
$ws.Range("H132").Value = 6746
$ws.Range("I132").Value = 3992.5
$ws.Range("K132").Value = 11977.5
$ws.Range("M132").Value = -9447.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6568.6924
$ws.Range("I46").Value = 5035.909
$ws.Range("K46").Value = 5035.909
$ws.Range("M46").Value = -4847.909
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
$ws.Range("H107").Value = 76924340
$ws.Range("I107").Value = 90910450
$ws.Range("J107").Value = 711.5
$ws.Range("K107").Value = 272731350
$ws.Range("L107").Value = 2134.5
$ws.Range("M107").Value = -272729430
$ws.Range("N107").Value = -5974.5
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null
$ws.Range("H123").Value = 65000
$ws.Range("J123").Value = 65000
$ws.Range("L123").Value = 65000
$ws.Range("N123").Value = -74800
$ws.Range("H132").Value = 36110210
$ws.Range("I132").Value = 66676310
$ws.Range("K132").Value = 200028930
$ws.Range("M132").Value = -200026400
